$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1864406779661017
$ws.Range("C2").Value = 0.576271186440678
$ws.Range("J2").Value = 0.01016949152542373
$ws.Range("P2").Value = 0.1423728813559322
$ws.Range("S2").Value = 0.0847457627118644
$ws.Range("C3").Value = 0.02840909090909091
$ws.Range("J3").Value = 0.03409090909090909
$ws.Range("P3").Value = 0.7613636363636364
$ws.Range("S3").Value = 0.1761363636363636
$ws.Range("J4").Value = 0.1081081081081081
$ws.Range("P4").Value = 0.7027027027027027
$ws.Range("S4").Value = 0.1891891891891892
$ws.Range("B6").Value = 0.06217616580310881
$ws.Range("D6").Value = 0.02590673575129534
$ws.Range("F6").Value = 0.05699481865284974
$ws.Range("J6").Value = 0.3160621761658031
$ws.Range("O6").Value = 0.0155440414507772
$ws.Range("Q6").Value = 0.1450777202072539
$ws.Range("R6").Value = 0.05699481865284974
$ws.Range("S6").Value = 0.3212435233160622
$ws.Range("B7").Value = 0.1313131313131313
$ws.Range("D7").Value = 0.0202020202020202
$ws.Range("F7").Value = 0.05555555555555555
$ws.Range("J7").Value = 0.1212121212121212
$ws.Range("O7").Value = 0.0303030303030303
$ws.Range("Q7").Value = 0.196969696969697
$ws.Range("R7").Value = 0.1262626262626263
$ws.Range("S7").Value = 0.3181818181818182
$ws.Range("B8").Value = 0.1002710027100271
$ws.Range("D8").Value = 0.02439024390243903
$ws.Range("E8").Value = 0.002710027100271003
$ws.Range("F8").Value = 0.06233062330623306
$ws.Range("J8").Value = 0.1138211382113821
$ws.Range("O8").Value = 0.01626016260162602
$ws.Range("Q8").Value = 0.1815718157181572
$ws.Range("R8").Value = 0.1002710027100271
$ws.Range("S8").Value = 0.3983739837398374
$ws.Range("B9").Value = 0.09090909090909091
$ws.Range("D9").Value = 0.01136363636363636
$ws.Range("E9").Value = 0.005681818181818182
$ws.Range("F9").Value = 0.09090909090909091
$ws.Range("J9").Value = 0.125
$ws.Range("O9").Value = 0.01136363636363636
$ws.Range("Q9").Value = 0.1534090909090909
$ws.Range("R9").Value = 0.1079545454545455
$ws.Range("S9").Value = 0.4034090909090909
$ws.Range("B10").Value = 0.1287744227353464
$ws.Range("D10").Value = 0.01509769094138544
$ws.Range("E10").Value = 0.0008880994671403197
$ws.Range("F10").Value = 0.0630550621669627
$ws.Range("J10").Value = 0.1145648312611012
$ws.Range("O10").Value = 0.01154529307282416
$ws.Range("Q10").Value = 0.2007104795737123
$ws.Range("R10").Value = 0.1012433392539965
$ws.Range("S10").Value = 0.3641207815275311
$ws.Range("G11").Value = 0.147887323943662
$ws.Range("J11").Value = 0.06690140845070422
$ws.Range("K11").Value = 0.176056338028169
$ws.Range("L11").Value = 0.5950704225352113
$ws.Range("S11").Value = 0.01408450704225352
$ws.Range("G12").Value = 0.7529411764705882
$ws.Range("J12").Value = 0.2
$ws.Range("K12").Value = 0.01176470588235294
$ws.Range("L12").Value = 0.01176470588235294
$ws.Range("S12").Value = 0.02352941176470588
$ws.Range("G13").Value = 0.7291666666666666
$ws.Range("J13").Value = 0.2708333333333333
$ws.Range("F15").Value = 0.02127659574468085
$ws.Range("H15").Value = 0.1436170212765958
$ws.Range("I15").Value = 0.09042553191489362
$ws.Range("J15").Value = 0.3829787234042553
$ws.Range("K15").Value = 0.06382978723404255
$ws.Range("M15").Value = 0.01595744680851064
$ws.Range("O15").Value = 0.03723404255319149
$ws.Range("S15").Value = 0.2446808510638298
$ws.Range("F16").Value = 0.03092783505154639
$ws.Range("H16").Value = 0.1494845360824742
$ws.Range("I16").Value = 0.07216494845360824
$ws.Range("J16").Value = 0.4484536082474227
$ws.Range("K16").Value = 0.1185567010309278
$ws.Range("M16").Value = 0.02061855670103093
$ws.Range("O16").Value = 0.04639175257731959
$ws.Range("S16").Value = 0.1134020618556701
$ws.Range("F17").Value = 0.01587301587301587
$ws.Range("H17").Value = 0.1402116402116402
$ws.Range("I17").Value = 0.09259259259259259
$ws.Range("J17").Value = 0.4259259259259259
$ws.Range("K17").Value = 0.1296296296296296
$ws.Range("M17").Value = 0.01851851851851852
$ws.Range("O17").Value = 0.06084656084656084
$ws.Range("S17").Value = 0.1164021164021164
$ws.Range("F18").Value = 0.01951219512195122
$ws.Range("H18").Value = 0.1804878048780488
$ws.Range("I18").Value = 0.08780487804878048
$ws.Range("J18").Value = 0.424390243902439
$ws.Range("K18").Value = 0.1512195121951219
$ws.Range("M18").Value = 0.01951219512195122
$ws.Range("O18").Value = 0.03902439024390244
$ws.Range("S18").Value = 0.07804878048780488
$ws.Range("F19").Value = 0.02296650717703349
$ws.Range("H19").Value = 0.215311004784689
$ws.Range("I19").Value = 0.08899521531100478
$ws.Range("J19").Value = 0.3559808612440191
$ws.Range("K19").Value = 0.108133971291866
$ws.Range("M19").Value = 0.03062200956937799
$ws.Range("O19").Value = 0.07942583732057416
$ws.Range("S19").Value = 0.09856459330143541
